# Keystone Advanced Debug.pptx - apply commit changes
#  - Add a presentation-level custom tag (ARTICULATE_PROJECT_OPEN=0)
#  - Update the cached "last printed" date field in the notes master from 2/24/2012 to 3/6/2012
#  - Remove the unused Date/Footer/SlideNumber placeholder shapes from the "Title Slide" layout
#  - Resize/reposition the "Rectangle 9" label on the slide master and retitle it
#    from "CI Training" to "Multicore Training"

$p = $ppt.ActivePresentation

# 1) New presentation-level Articulate tag -> creates <p:custDataLst><p:tags .../></p:custDataLst>
#    directly under <p:presentation>.
$p.Tags.Add("ARTICULATE_PROJECT_OPEN", "0")

# 2) Notes master date placeholder: cached date text 2/24/2012 -> 3/6/2012
$notesMaster = $p.NotesMaster
$dateAndTime = $notesMaster.HeadersFooters.DateAndTime
$dateAndTime.Text = "3/6/2012"

# 3) Remove the Date/Footer/Slide Number placeholders from the 3rd slide layout
#    ("Title Slide", used by the slide master) - delete from the end so indices
#    of the remaining shapes stay valid while iterating.
$master = $p.SlideMaster
$titleLayout = $master.CustomLayouts.Item(3)
for ($i = $titleLayout.Shapes.Count; $i -ge 1; $i--) {
    $shape = $titleLayout.Shapes.Item($i)
    if ($shape.Name -eq "Date Placeholder 3" -or $shape.Name -eq "Footer Placeholder 4" -or $shape.Name -eq "Slide Number Placeholder 5") {
        $shape.Delete()
    }
}

# 4) Slide master "Rectangle 9" label: reposition/resize and rename text
$rect = $master.Shapes.Item("Rectangle 9")
$rect.Left = 7425393 / 12700
$rect.Width = 1357103 / 12700
$rect.Top = 6498264 / 12700
$rect.Height = 276999 / 12700
$rect.TextFrame.TextRange.Text = "Multicore Training"
